# "Iteration upto Water Coil"
# 1) Rename the existing sheet "Sheet2" -> "Casing" and tweak its box-dimension inputs.
# 2) Add a new "Coil" sheet (water-coil sizing calc) after it and make it the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Casing sheet: rename + update inputs
# ---------------------------------------------------------------------------
$casing = $wb.Worksheets.Item(1)
$casing.Name = "Casing"

$casing.Range("A2").Value = 0.6
$casing.Range("C2").Value = 1000
$casing.Range("D2").Value = 500
$casing.Range("E2").Value = 500

$casing.Range("P5").Select()

# ---------------------------------------------------------------------------
# 2. New Coil sheet
# ---------------------------------------------------------------------------
$coil = $wb.Worksheets.Add()
$coil.Name = "Coil"
$coil.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$coil = $wb.Worksheets.Item("Coil")

# Header row (A1:AG1)
$headers = @(
    "fin_height","fin_length","rd","tubes","fpi","fin_qty","tube_qty",
    "casing_qty","header_qty","drainpan_sheet_qty","U_bend","copper_st_stub",
    "fin_qty_f1","fin_qty_f2","casing_qty_f1","casing_qty_f2","casing_qty_f3",
    "casing_qty_f4","casing_qty_f5","casing_qty_f6","tubes_f1","header_qty_f1",
    "header_qty_f2","header_qty_f3","drainpan_qty_f1","drainpan_qty_f2",
    "drainpan_qty_f3","U_bend_f1","U_bend_f2","U_bend_f3","copper_st_stub_f1",
    "copper_st_stub_f2","meter_conv"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $coil.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - raw inputs
$coil.Range("A2").Value = 1000
$coil.Range("B2").Value = 1200
$coil.Range("C2").Value = 6
$coil.Range("E2").Value = 12
$coil.Range("M2").Value = 0.45
$coil.Range("N2").Value = 10.76
$coil.Range("O2").Value = 2
$coil.Range("P2").Value = 1
$coil.Range("Q2").Value = 40
$coil.Range("R2").Value = 7.81
$coil.Range("S2").Value = 1.6
$coil.Range("T2").Value = 1.15
$coil.Range("U2").Value = 31.75
$coil.Range("V2").Value = 2
$coil.Range("W2").Value = 350
$coil.Range("X2").Value = 2
$coil.Range("Y2").Value = 0.6
$coil.Range("Z2").Value = 7.81
$coil.Range("AA2").Value = 1
$coil.Range("AB2").Value = 2
$coil.Range("AC2").Value = 2
$coil.Range("AD2").Value = 2
$coil.Range("AE2").Value = 2
$coil.Range("AF2").Value = 0.015
$coil.Range("AG2").Value = 1000

# Row 2 - formulas
$coil.Range("D2").Formula = "=+A2/U2"
$coil.Range("F2").Formula = "=(A2/AG2)*(B2/AG2)*C2*M2*N2"
$coil.Range("G2").Formula = "=(A2/AG2)*(B2/AG2)*C2*M2*N2"
$coil.Range("H2").Formula = "=((A2+B2)*O2/AG2)*((C2+P2)*Q2/AG2)*R2*S2*T2"
$coil.Range("I2").Formula = "=IF(AND(A2>0,B2>0),(A2*V2+W2)*X2/AG2,0)"
$coil.Range("J2").Formula = "=IF(AND(A2>0,B2>0),(Y2*B2/AG2*Z2*AA2),0)"
$coil.Range("K2").Formula = "=+A2/U2*(C2*AB2-AC2)/AD2"
$coil.Range("L2").Formula = "=+D2*AE2*AF2"

# D2 / K2 carry an integer ("0") number format
$coil.Range("D2").NumberFormat = "0"
$coil.Range("K2").NumberFormat = "0"

# Column widths (approximate auto-fit from the authored sheet)
$coil.Range("A1:B1").ColumnWidth = 10.140625
$coil.Range("C1").ColumnWidth = 2.85546875
$coil.Range("D1").ColumnWidth = 6
$coil.Range("E1").ColumnWidth = 3.42578125
$coil.Range("F1").ColumnWidth = 8
$coil.Range("H1").ColumnWidth = 11
$coil.Range("I1").ColumnWidth = 11.140625
$coil.Range("J1").ColumnWidth = 19.42578125
$coil.Range("K1").ColumnWidth = 7.85546875
$coil.Range("L1").ColumnWidth = 14.7109375
$coil.Range("M1:Q1").ColumnWidth = 10
$coil.Range("R1:T1").ColumnWidth = 13.140625
$coil.Range("U1").ColumnWidth = 8.7109375
$coil.Range("V1:X1").ColumnWidth = 14
$coil.Range("Y1:AA1").ColumnWidth = 15.5703125
$coil.Range("AB1:AD1").ColumnWidth = 10.5703125
$coil.Range("AE1:AF1").ColumnWidth = 17.5703125

$coil.Range("Z10").Select()

# Coil is the sheet the author left active
$coil.Activate()
